$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.605.61"
$ws.Range("E2").Value = "  -12.69%  "

$ws.Range("D3").Value = "2.341.88"
$ws.Range("E3").Value = "  -18.86%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "439.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -15.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "122.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.65%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.37%  "

$ws.Range("E8").Value = "  -13.75%  "

$ws.Range("D9").Value = "2.349.12"
$ws.Range("E9").Value = "  -18.81%  "

$ws.Range("E10").Value = "  -12.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0895"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -16.31%  "

$ws.Range("E12").Value = "  -13.79%  "

$ws.Range("E13").Value = "  -5.28%  "

$ws.Range("D14").Value = "52.596.56"
$ws.Range("E14").Value = "  -12.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.13"
$ws.Range("D15").Style = "Normal"

$ws.Range("E16").Value = "  -14.35%  "

$ws.Range("D17").Value = "2.357.46"
$ws.Range("E17").Value = "  -18.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -18.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "303.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -14.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -21.56%  "

$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("E22").Value = "  -1.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -20.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "54.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -15.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.154"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -13.81%  "

$ws.Range("E26").Value = "  -17.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.38%  "

$ws.Range("D29").Value = "0.0₃0686"
$ws.Range("E29").Value = "  -17.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "144.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "17.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -12.26%  "

$ws.Range("E32").Value = "  -19.39%  "

$ws.Range("E33").Value = "  -12.72%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.836"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -15.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -17.66%  "

$ws.Range("E36").Value = "  -15.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.994"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "32.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -14.21%  "

$ws.Range("E39").Value = "  -13.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.15"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0513"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -11.65%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -15.86%  "

$ws.Range("D43").Value = "1.935.46"
$ws.Range("E43").Value = "  -15.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.534"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -17.19%  "

$ws.Range("E45").Value = "  -11.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0834"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -15.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -21.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.62%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -12.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "15.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -15.31%  "
